# Sprint Planning Meeting (Sprint 2) - content / layout refresh
# (re-save normalised the runs, relocated the stray _GoBack bookmark and
#  switched the page template from US Letter to A4, plus a couple of
#  Normal-style tweaks that came along with the new template.)

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Collapse the fragmented runs of the "second sprint" paragraph into
#    a single run. Doing the replacement through Find/Execute (instead
#    of touching .Text) is what makes Word re-flow the paragraph down
#    to one <w:r>, exactly like a normal editing pass would.
# ---------------------------------------------------------------------
$para2 = $d.Paragraphs(2).Range
$para2Text = "In the second sprint, we will meet with the client and collect more information, feedback and user story to ensure to meet their needs. In this sprint, we will try to focus on the task and feedback we have received from the client. We decide to rush up a bit since we have existing interface. "
$para2.Find.Execute($para2Text, $false, $false, $false, $false, $false, $true, 1, $false, $para2Text, 2)

# ---------------------------------------------------------------------
# 2) Same treatment for the "first week / second week / third week"
#    paragraph. This also swallows the _GoBack bookmark that used to
#    sit in the middle of the text (it gets recreated at the end of the
#    document in step 3, matching where Word leaves it after an edit).
# ---------------------------------------------------------------------
$para4 = $d.Paragraphs(4).Range
$para4Text = "In the first week we will be cover the report first where we ensure it can retrieve the data and ensure that other module are proceed as usual with the new user story.  In the second week, we try ensure bus schedule module, staff module, route module is completed so that 3 developer can process to do the bus module as soon as possible if there is no problem and will have feedback from client if necessary. In the third week, we must complete all user story in time and get the feedback from client and improve the system based on their story. Lastly, we will hope to complete it on time and successfully deliver the system to the client"
$para4.Find.Execute($para4Text, $false, $false, $false, $false, $false, $true, 1, $false, $para4Text, 2)

# ---------------------------------------------------------------------
# 3) Recreate the _GoBack bookmark as its own empty paragraph at the
#    very end of the document (this is where Word's "last edit"
#    position bookmark ends up after the text above was retyped).
# ---------------------------------------------------------------------
$endRange = $d.Range($d.Content.End - 1, $d.Content.End - 1)
$endRange.Text = "`r"

$lastPara = $d.Paragraphs($d.Paragraphs.Count)
$lastPara.Alignment = 0

$d.Bookmarks.Add("_GoBack", $lastPara.Range)

# ---------------------------------------------------------------------
# 4) Page setup: switch from the US-Letter template to the A4 template
#    (page size, margins, header/footer distance and column spacing).
# ---------------------------------------------------------------------
$ps = $d.PageSetup
$ps.PageWidth = 595.3
$ps.PageHeight = 841.9
$ps.TopMargin = 72
$ps.RightMargin = 90
$ps.BottomMargin = 72
$ps.LeftMargin = 90
$ps.HeaderDistance = 42.55
$ps.FooterDistance = 49.6
$ps.Gutter = 0
$ps.TextColumns.Spacing = 21.25

# ---------------------------------------------------------------------
# 5) Normal style picked up a couple of properties from the new
#    template as well (kerning off, explicit run size, paragraph
#    spacing/line rule).
# ---------------------------------------------------------------------
$normal = $d.Styles("Normal")
$normal.Font.Size = 11
$normal.Font.Kerning = 0
$normal.ParagraphFormat.SpaceAfter = 8
$normal.ParagraphFormat.LineSpacingRule = 5
$normal.ParagraphFormat.LineSpacing = 12.8
